# The presentation had a slide (in position 22, sldId 275 — the
# "/ Evaluations and Demo / Future Improvements" slide with the POS
# Features / Pre-Processing callouts) removed. The two slides that
# followed it ("Live Demo" and "Thanks for Listening!") simply shift up
# by one position as a natural consequence of the deletion.

$p = $ppt.ActivePresentation

$s = $p.Slides.Item(22)
$s.Delete()
